$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 currently holds phone "09876543" as text, with total_points = 0.
# Duplicate that row down into a new row 13 so the original text value is
# preserved there, then convert row 12's phone value into a genuine number
# (9876543) while keeping its total_points at 0.

$ws.Range("A12:C12").Copy()
$ws.Range("A13:C13").PasteSpecial()

$ws.Cells.Item(12, 1).Value = 9876543
$ws.Cells.Item(12, 3).Value = 0
